$d = $word.ActiveDocument

$d.Content.Find.Execute(".park}", $true, $false, $false, $false, $false, $true, 1, $false, ".park.name}", 2)
$d.Content.Find.Execute("{d.area}", $true, $false, $false, $false, $false, $true, 1, $false, "{d.parkAreasFormatted}", 2)
